# Update "想去人数" (want-to-go count) figures on gh-pages data refresh.
# Mirrors commit: "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$wsExpo  = $wb.Worksheets.Item("展览")
$wsShow  = $wb.Worksheets.Item("演出")
$wsAll   = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions) sheet
$wsExpo.Range("F12").Value = 4540
$wsExpo.Range("F13").Value = 4540
$wsExpo.Range("F20").Value = 7103
$wsExpo.Range("F39").Value = 90
$wsExpo.Range("F48").Value = 2096
$wsExpo.Range("F50").Value = 1054

# 演出 (Shows) sheet
$wsShow.Range("F3").Value = 218
$wsShow.Range("F5").Value = 37

# 全部类型 (All types) sheet
$wsAll.Range("F8").Value = 218
$wsAll.Range("F10").Value = 37
$wsAll.Range("F12").Value = 4540
$wsAll.Range("F13").Value = 4540
$wsAll.Range("F20").Value = 7103
$wsAll.Range("F41").Value = 90
$wsAll.Range("F49").Value = 2096
$wsAll.Range("F51").Value = 1054
